$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item(5)

# Header row (row 1)
$ws5.Range("B1").Value = "company"
$ws5.Range("C1").Value = "name"
$ws5.Range("D1").Value = "owner"
$ws5.Range("E1").Value = "property_category"
$ws5.Range("F1").Value = "category"
$ws5.Range("G1").Value = "date"
$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("J1").Value = "source_file"
$ws5.Range("K1").Value = "index"

# Row 2
$ws5.Range("B2").Value = "新光人壽"
$ws5.Range("C2").Value = "新長安終身壽險（100000元）"
$ws5.Range("D2").Value = "鄭天財"
$ws5.Range("E2").Value = "insurance"
$ws5.Range("F2").Value = "normal"
$ws5.Range("G2").Value = "2012-04-30"
$ws5.Range("H2").Value = "鄭天財"
$ws5.Range("I2").Value = 1763
$ws5.Range("J2").Value = "tmp1c9c1"
$ws5.Range("K2").Value = 106

# Row 3
$ws5.Range("B3").Value = "新光人壽"
$ws5.Range("C3").Value = "百年長青終身壽險(1000000元）"
$ws5.Range("D3").Value = "鄭天財"
$ws5.Range("E3").Value = "insurance"
$ws5.Range("F3").Value = "normal"
$ws5.Range("G3").Value = "2012-04-30"
$ws5.Range("H3").Value = "鄭天財"
$ws5.Range("I3").Value = 1763
$ws5.Range("J3").Value = "tmp1c9c1"
$ws5.Range("K3").Value = 108

# Row 4
$ws5.Range("B4").Value = "台灣人壽"
$ws5.Range("C4").Value = "新長榮還本終身壽險(600000元）"
$ws5.Range("D4").Value = "王慧玲"
$ws5.Range("E4").Value = "insurance"
$ws5.Range("F4").Value = "normal"
$ws5.Range("G4").Value = "2012-04-30"
$ws5.Range("H4").Value = "鄭天財"
$ws5.Range("I4").Value = 1763
$ws5.Range("J4").Value = "tmp1c9c1"
$ws5.Range("K4").Value = 109

# Row 5
$ws5.Range("B5").Value = "安泰人壽"
$ws5.Range("C5").Value = "雙星報喜還本終身壽險"
$ws5.Range("D5").Value = "王慧玲"
$ws5.Range("E5").Value = "insurance"
$ws5.Range("F5").Value = "normal"
$ws5.Range("G5").Value = "2012-04-30"
$ws5.Range("H5").Value = "鄭天財"
$ws5.Range("I5").Value = 1763
$ws5.Range("J5").Value = "tmp1c9c1"
$ws5.Range("K5").Value = 110

# Row 6
$ws5.Range("B6").Value = "安泰人壽"
$ws5.Range("C6").Value = "雙星報喜還本終身壽險"
$ws5.Range("D6").Value = "王慧玲"
$ws5.Range("E6").Value = "insurance"
$ws5.Range("F6").Value = "normal"
$ws5.Range("G6").Value = "2012-04-30"
$ws5.Range("H6").Value = "鄭天財"
$ws5.Range("I6").Value = 1763
$ws5.Range("J6").Value = "tmp1c9c1"
$ws5.Range("K6").Value = 111

$ws6 = $wb.Worksheets.Item(6)

# Header row (row 1)
$ws6.Range("B1").Value = "species"
$ws6.Range("C1").Value = "owner"
$ws6.Range("D1").Value = "debtor"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

# Row 2
$ws6.Range("B2").Value = "般借款"
$ws6.Range("C2").Value = "鄭天財"
$ws6.Range("D2").Value = "陳俊民南投縣南投市中興路"
$ws6.Range("E2").Value = 1389488
$ws6.Range("F2").Value = "95年02月25日"
$ws6.Range("G2").Value = "南投地方法院支付命令"
$ws6.Range("H2").Value = "claim"
$ws6.Range("I2").Value = "normal"
$ws6.Range("J2").Value = "2012-04-30"
$ws6.Range("K2").Value = "鄭天財"
$ws6.Range("L2").Value = 1763
$ws6.Range("M2").Value = "tmp1c9c1"
$ws6.Range("N2").Value = 116

$ws7 = $wb.Worksheets.Item(7)

# Header row (row 1)
$ws7.Range("B1").Value = "species"
$ws7.Range("C1").Value = "debtor"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"

# Row 2
$ws7.Range("B2").Value = "房屋貸款"
$ws7.Range("C2").Value = "鄭天財"
$ws7.Range("D2").Value = "台灣土地銀行士林分行臺北市土林區中山北路"
$ws7.Range("E2").Value = 991281
$ws7.Range("F2").Value = "95年06月26日"
$ws7.Range("G2").Value = "購置房屋"
$ws7.Range("H2").Value = "debt"
$ws7.Range("I2").Value = "normal"
$ws7.Range("J2").Value = "2012-04-30"
$ws7.Range("K2").Value = "鄭天財"
$ws7.Range("L2").Value = 1763
$ws7.Range("M2").Value = "tmp1c9c1"
$ws7.Range("N2").Value = 121

